# Generate Report for Archive
#
# The "b01085d3-d71c-45fa-8065-8d797978bd65" file entry moves from the last
# data row (row 9) to become the first of the "alphabetically later" group
# (row 6), right before the "e2aa29d7..." row. The rows that used to sit at
# 6, 7 and 8 each shift down by one (to 7, 8, 9). This same rotation happens
# identically on all three worksheets: Overview, zh-cn and de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "b01085d3-d71c-45fa-8065-8d797978bd65.md"
$wsOverview.Range("B6").Value = "e2e\b01085d3-d71c-45fa-8065-8d797978bd65.md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-21 12:49:58"

$wsOverview.Range("A7").Value = "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md"
$wsOverview.Range("B7").Value = "e2e\e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md"
$wsOverview.Range("E7").Value = "In Translation"
$wsOverview.Range("F7").Value = "In Translation"
$wsOverview.Range("G7").Value = "2016-08-21 12:47:00"

$wsOverview.Range("A8").Value = "3a3ae932-91a2-44b2-b732-dfd2131fb523.md"
$wsOverview.Range("B8").Value = "e2e\3a3ae932-91a2-44b2-b732-dfd2131fb523.md"
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = "2016-08-21 12:45:53"

$wsOverview.Range("A9").Value = "4e4d33e6-0c2a-4994-854e-75a7039c8d10.md"
$wsOverview.Range("B9").Value = "e2e\4e4d33e6-0c2a-4994-854e-75a7039c8d10.md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-21 12:47:58"

# Rebuild the hyperlinks collection in the new row order so the display
# text tracks the new file name shown in column B, preserving the
# relationship order (rId2..rId9).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f7e8ef10add7166f54a2560fdb144972762ee04e/e2e/f02d3662-775a-4d14-b928-c22c4c2a93eb.md", "", "", "e2e\f02d3662-775a-4d14-b928-c22c4c2a93eb.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/0910a0c4-f35c-4f0b-b270-9213140b88f6.md", "", "", "e2e\0910a0c4-f35c-4f0b-b270-9213140b88f6.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3105897da6e243f7864e9573006ee92fa59cc9/e2e/0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md", "", "", "e2e\0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0c0949ce25096d6d474b01abe263efaf3022a6/e2e/6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md", "", "", "e2e\6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d94689885726640f0d7b9cd7d0e40f05a95b0cf/e2e/b01085d3-d71c-45fa-8065-8d797978bd65.md", "", "", "e2e\b01085d3-d71c-45fa-8065-8d797978bd65.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md", "", "", "e2e\e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8e5feb1f7649373cece8ad20c883fb123bd5cb1/e2e/3a3ae932-91a2-44b2-b732-dfd2131fb523.md", "", "", "e2e\3a3ae932-91a2-44b2-b732-dfd2131fb523.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e8fedd5305e8b5099a3cd0826c50f7aeab8b119/e2e/4e4d33e6-0c2a-4994-854e-75a7039c8d10.md", "", "", "e2e\4e4d33e6-0c2a-4994-854e-75a7039c8d10.md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "b01085d3-d71c-45fa-8065-8d797978bd65.md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("G6").Value = "b01085d3-d71c-45fa-8065-8d797978bd65.c154ab9eafb3d2750d2d14205a23953ec80ce3be.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-21 12:49:54"

$wsZhCn.Range("A7").Value = "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md"
$wsZhCn.Range("C7").Value = "In Translation"
$wsZhCn.Range("G7").Value = "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.09fc332fcc3d5200af76ac5db7db85e8b631eb8d.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-21 12:46:55"

$wsZhCn.Range("A8").Value = "3a3ae932-91a2-44b2-b732-dfd2131fb523.md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("G8").Value = "3a3ae932-91a2-44b2-b732-dfd2131fb523.0093d287a44e4bfdd9f66c0707e42d0b082b8957.zh-cn.xlf"
$wsZhCn.Range("H8").Value = "2016-08-21 12:45:49"

$wsZhCn.Range("A9").Value = "4e4d33e6-0c2a-4994-854e-75a7039c8d10.md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("G9").Value = "4e4d33e6-0c2a-4994-854e-75a7039c8d10.50007debd7244bb22ca916e7fc63dedd44fe3af3.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-21 12:47:54"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f7e8ef10add7166f54a2560fdb144972762ee04e/e2e/f02d3662-775a-4d14-b928-c22c4c2a93eb.md", "", "", "f02d3662-775a-4d14-b928-c22c4c2a93eb.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1e1d96e557cbaa3134e024b65d33886d0ec63736/e2e/f02d3662-775a-4d14-b928-c22c4c2a93eb.md", "", "", "f02d3662-775a-4d14-b928-c22c4c2a93eb.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/0910a0c4-f35c-4f0b-b270-9213140b88f6.md", "", "", "0910a0c4-f35c-4f0b-b270-9213140b88f6.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3105897da6e243f7864e9573006ee92fa59cc9/e2e/0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md", "", "", "0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/43663c22ef1fd9bbb90194b4aec4fb6e29b6d012/e2e/0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md", "", "", "0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0c0949ce25096d6d474b01abe263efaf3022a6/e2e/6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md", "", "", "6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d94689885726640f0d7b9cd7d0e40f05a95b0cf/e2e/b01085d3-d71c-45fa-8065-8d797978bd65.md", "", "", "b01085d3-d71c-45fa-8065-8d797978bd65.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md", "", "", "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8e5feb1f7649373cece8ad20c883fb123bd5cb1/e2e/3a3ae932-91a2-44b2-b732-dfd2131fb523.md", "", "", "3a3ae932-91a2-44b2-b732-dfd2131fb523.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e8fedd5305e8b5099a3cd0826c50f7aeab8b119/e2e/4e4d33e6-0c2a-4994-854e-75a7039c8d10.md", "", "", "4e4d33e6-0c2a-4994-854e-75a7039c8d10.md")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "b01085d3-d71c-45fa-8065-8d797978bd65.md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("G6").Value = "b01085d3-d71c-45fa-8065-8d797978bd65.c154ab9eafb3d2750d2d14205a23953ec80ce3be.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-21 12:49:58"

$wsDeDe.Range("A7").Value = "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md"
$wsDeDe.Range("C7").Value = "In Translation"
$wsDeDe.Range("G7").Value = "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.09fc332fcc3d5200af76ac5db7db85e8b631eb8d.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-21 12:47:00"

$wsDeDe.Range("A8").Value = "3a3ae932-91a2-44b2-b732-dfd2131fb523.md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("G8").Value = "3a3ae932-91a2-44b2-b732-dfd2131fb523.0093d287a44e4bfdd9f66c0707e42d0b082b8957.de-de.xlf"
$wsDeDe.Range("H8").Value = "2016-08-21 12:45:53"

$wsDeDe.Range("A9").Value = "4e4d33e6-0c2a-4994-854e-75a7039c8d10.md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("G9").Value = "4e4d33e6-0c2a-4994-854e-75a7039c8d10.50007debd7244bb22ca916e7fc63dedd44fe3af3.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-21 12:47:58"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f7e8ef10add7166f54a2560fdb144972762ee04e/e2e/f02d3662-775a-4d14-b928-c22c4c2a93eb.md", "", "", "f02d3662-775a-4d14-b928-c22c4c2a93eb.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c4ebbbab6ad0d4bb1952df656b5d0f1520147427/e2e/f02d3662-775a-4d14-b928-c22c4c2a93eb.md", "", "", "f02d3662-775a-4d14-b928-c22c4c2a93eb.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/0910a0c4-f35c-4f0b-b270-9213140b88f6.md", "", "", "0910a0c4-f35c-4f0b-b270-9213140b88f6.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3105897da6e243f7864e9573006ee92fa59cc9/e2e/0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md", "", "", "0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/71b5b3deb6140bf10150645b4560a3192f4715be/e2e/0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md", "", "", "0d84ec85-f7c5-4ee7-806f-69da21a8e48a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c0c0949ce25096d6d474b01abe263efaf3022a6/e2e/6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md", "", "", "6d62b4b1-0876-4dc7-acd6-f46d2abd144a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d94689885726640f0d7b9cd7d0e40f05a95b0cf/e2e/b01085d3-d71c-45fa-8065-8d797978bd65.md", "", "", "b01085d3-d71c-45fa-8065-8d797978bd65.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec966f340a0539fc47ef1d32dc2da6dd091bff/e2e/e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md", "", "", "e2aa29d7-b4c2-48bf-b6ee-d69b02e00f8c.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8e5feb1f7649373cece8ad20c883fb123bd5cb1/e2e/3a3ae932-91a2-44b2-b732-dfd2131fb523.md", "", "", "3a3ae932-91a2-44b2-b732-dfd2131fb523.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e8fedd5305e8b5099a3cd0826c50f7aeab8b119/e2e/4e4d33e6-0c2a-4994-854e-75a7039c8d10.md", "", "", "4e4d33e6-0c2a-4994-854e-75a7039c8d10.md")
